$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The finder's test resources moved under a SampleImageSetByExcelFileFinder
# sub-folder, so every recorded image path in this template must be updated
# to include that new path segment.
$null = $ws.Cells.Replace("excel/OneCamera/", "excel/SampleImageSetByExcelFileFinder/OneCamera/")

# Leave the selection where it ended up after editing the sheet.
$ws.Range("J28").Select() | Out-Null
